$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 1241-1242; this pushes the former rows 1241:1305
# down to 1243:1307 (dimension grows from A1:R1305 to A1:R1307), matching
# every downstream row's new content.
$ws.Rows("1241:1242").Insert()

# Row 1241 - "Primera" quality entry for the new week (boilerplate columns
# A,B,C,E,F,G,H,N,O,Q,R copied from the repeating pattern used by every
# other row in this table).
$ws.Cells.Item(1241, 1).Value = 8
$ws.Cells.Item(1241, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(1241, 3).Value = "Coquimbo"
$ws.Cells.Item(1241, 4).Value = 45267
$ws.Cells.Item(1241, 5).Value = 4
$ws.Cells.Item(1241, 6).Value = 100112023
$ws.Cells.Item(1241, 7).Value = "Brócoli"
$ws.Cells.Item(1241, 8).Value = "Sin especificar"
$ws.Cells.Item(1241, 9).Value = "Primera"
$ws.Cells.Item(1241, 10).Value = 2000
$ws.Cells.Item(1241, 11).Value = 800
$ws.Cells.Item(1241, 12).Value = 900
$ws.Cells.Item(1241, 13).Value = 850
$ws.Cells.Item(1241, 14).Value = "$/unidad"
$ws.Cells.Item(1241, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(1241, 16).Value = 850
$ws.Cells.Item(1241, 17).Value = 1
$ws.Cells.Item(1241, 18).Value = "Hortaliza"

# Row 1242 - "Segunda" quality entry for the same new week.
$ws.Cells.Item(1242, 1).Value = 8
$ws.Cells.Item(1242, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(1242, 3).Value = "Coquimbo"
$ws.Cells.Item(1242, 4).Value = 45267
$ws.Cells.Item(1242, 5).Value = 4
$ws.Cells.Item(1242, 6).Value = 100112023
$ws.Cells.Item(1242, 7).Value = "Brócoli"
$ws.Cells.Item(1242, 8).Value = "Sin especificar"
$ws.Cells.Item(1242, 9).Value = "Segunda"
$ws.Cells.Item(1242, 10).Value = 1160
$ws.Cells.Item(1242, 11).Value = 600
$ws.Cells.Item(1242, 12).Value = 700
$ws.Cells.Item(1242, 13).Value = 650
$ws.Cells.Item(1242, 14).Value = "$/unidad"
$ws.Cells.Item(1242, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(1242, 16).Value = 650
$ws.Cells.Item(1242, 17).Value = 1
$ws.Cells.Item(1242, 18).Value = "Hortaliza"
